$d = $word.ActiveDocument

# New (partially translated) campaign-dates sentence that replaces the old
# "Informace v této příručce ... souhvězdí Persea." paragraph content
# everywhere it occurs in the document body.
$newText = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od Pegasus: 8.-17. října, 7.-16. listopadu,"

# Collect every paragraph whose text still holds the old sentence. We grab
# the Paragraph objects themselves (not raw numbers) so that, as earlier
# paragraphs are edited and the document shrinks, the Start/End of the
# paragraphs we haven't processed yet stay correct automatically.
$targets = @()
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Informace v této příručce*") {
        $targets += $p
    }
}

foreach ($p in $targets) {
    $start = $p.Range.Start
    # Exclude the trailing paragraph mark (the last character of a
    # Paragraph.Range) so only the paragraph's run content is removed and
    # the paragraph itself (its pPr/paraId) survives.
    $end = $p.Range.End - 1

    $body = $d.Range($start, $end)
    $body.Delete()

    # Re-insert the replacement text into the now-empty paragraph. Because
    # this position carries no leftover run formatting, the newly created
    # run comes out plain (no rPr/rFonts/lang overrides), matching the
    # target markup.
    $insertion = $d.Range($start, $start)
    $insertion.InsertAfter($newText)
}
